$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 82.987681
$ws.Cells.Item(2, 8).Value = 248.963043
$ws.Cells.Item(2, 9).Value = 0.4489504115427952
$ws.Cells.Item(2, 10).Value = 0.4489504115427952
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 4.536281
$ws.Cells.Item(2, 14).Value = 13.608843
$ws.Cells.Item(2, 15).Value = 0.0729209419042404
$ws.Cells.Item(2, 16).Value = 0.07292094190424041
$ws.Cells.Item(2, 17).Value = 376.455440554361
$ws.Cells.Item(2, 18).Value = 3388.098964989249
$ws.Cells.Item(2, 19).Value = 0.03273788687799699
$ws.Cells.Item(2, 20).Value = 0.032737886877997
$ws.Cells.Item(3, 7).Value = 82.987681
$ws.Cells.Item(3, 8).Value = 248.963043
$ws.Cells.Item(3, 9).Value = 0.4489504115427952
$ws.Cells.Item(3, 10).Value = 0.4489504115427952
$ws.Cells.Item(3, 15).Value = 0.7154667412877611
$ws.Cells.Item(3, 16).Value = 0.7154667412877612
$ws.Cells.Item(3, 17).Value = 3693.607628480383
$ws.Cells.Item(3, 18).Value = 33242.46865632344
$ws.Cells.Item(3, 19).Value = 0.321209087946323
$ws.Cells.Item(3, 20).Value = 0.321209087946323
$ws.Cells.Item(4, 7).Value = 82.987681
$ws.Cells.Item(4, 8).Value = 248.963043
$ws.Cells.Item(4, 9).Value = 0.4489504115427952
$ws.Cells.Item(4, 10).Value = 0.4489504115427952
$ws.Cells.Item(4, 15).Value = 0.2116123168079984
$ws.Cells.Item(4, 16).Value = 0.2116123168079984
$ws.Cells.Item(4, 17).Value = 1092.451713738103
$ws.Cells.Item(4, 18).Value = 9832.065423642925
$ws.Cells.Item(4, 19).Value = 0.09500343671847526
$ws.Cells.Item(4, 20).Value = 0.09500343671847526
$ws.Cells.Item(5, 7).Value = 63.14058933333333
$ws.Cells.Item(5, 9).Value = 0.3415807409566563
$ws.Cells.Item(5, 10).Value = 0.3415807409566563
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 4.536281
$ws.Cells.Item(5, 14).Value = 13.608843
$ws.Cells.Item(5, 15).Value = 0.0729209419042404
$ws.Cells.Item(5, 16).Value = 0.07292094190424041
$ws.Cells.Item(5, 17).Value = 286.4234557216026
$ws.Cells.Item(5, 18).Value = 2577.811101494424
$ws.Cells.Item(5, 19).Value = 0.02490838936690772
$ws.Cells.Item(5, 20).Value = 0.02490838936690773
$ws.Cells.Item(6, 7).Value = 63.14058933333333
$ws.Cells.Item(6, 9).Value = 0.3415807409566563
$ws.Cells.Item(6, 10).Value = 0.3415807409566563
$ws.Cells.Item(6, 15).Value = 0.7154667412877611
$ws.Cells.Item(6, 16).Value = 0.7154667412877612
$ws.Cells.Item(6, 18).Value = 25292.29683927574
$ws.Cells.Item(6, 19).Value = 0.2443896596189178
$ws.Cells.Item(6, 20).Value = 0.2443896596189178
$ws.Cells.Item(7, 7).Value = 63.14058933333333
$ws.Cells.Item(7, 9).Value = 0.3415807409566563
$ws.Cells.Item(7, 10).Value = 0.3415807409566563
$ws.Cells.Item(7, 15).Value = 0.2116123168079984
$ws.Cells.Item(7, 16).Value = 0.2116123168079984
$ws.Cells.Item(7, 19).Value = 0.07228269197083079
$ws.Cells.Item(7, 20).Value = 0.0722826919708308
$ws.Cells.Item(8, 9).Value = 0.2094688475005485
$ws.Cells.Item(8, 10).Value = 0.2094688475005485
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 4.536281
$ws.Cells.Item(8, 14).Value = 13.608843
$ws.Cells.Item(8, 15).Value = 0.0729209419042404
$ws.Cells.Item(8, 16).Value = 0.07292094190424041
$ws.Cells.Item(8, 17).Value = 175.644537215702
$ws.Cells.Item(8, 18).Value = 1580.800834941318
$ws.Cells.Item(8, 19).Value = 0.01527466565933569
$ws.Cells.Item(8, 20).Value = 0.01527466565933569
$ws.Cells.Item(9, 9).Value = 0.2094688475005485
$ws.Cells.Item(9, 10).Value = 0.2094688475005485
$ws.Cells.Item(9, 15).Value = 0.7154667412877611
$ws.Cells.Item(9, 16).Value = 0.7154667412877612
$ws.Cells.Item(9, 17).Value = 1723.343409794979
$ws.Cells.Item(9, 19).Value = 0.1498679937225204
$ws.Cells.Item(9, 20).Value = 0.1498679937225204
$ws.Cells.Item(10, 9).Value = 0.2094688475005485
$ws.Cells.Item(10, 10).Value = 0.2094688475005485
$ws.Cells.Item(10, 15).Value = 0.2116123168079984
$ws.Cells.Item(10, 16).Value = 0.2116123168079984
$ws.Cells.Item(10, 17).Value = 509.7101941400188
$ws.Cells.Item(10, 18).Value = 4587.391747260169
$ws.Cells.Item(10, 19).Value = 0.04432618811869237
$ws.Cells.Item(10, 20).Value = 0.04432618811869238
